$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE: 2020-09-30 -> 2019-09-30
$ws.Range("H2").Value = "2019-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.3847
# DEDUCT_BASIC_EPS
$ws.Range("J2").Value = 0.3739
# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 322849048.67
# PARENT_NETPROFIT
$ws.Range("L2").Value = 51934970.67
# WEIGHTAVG_ROE
$ws.Range("M2").Value = 22.59

# YSTZ, SJLTZ, BPS, MGJYXJJE cleared out
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("Q2").ClearContents()

# XSMLL
$ws.Range("R2").Value = 40.3755977622

# ISNEW (force text so "0" is not stored as a number)
$ws.Range("AB2").Value = "'0"
# QDATE
$ws.Range("AC2").Value = "2019Q3"
# DATATYPE
$ws.Range("AD2").Value = "2019年 三季报"
# DATAYEAR (force text so "2019" is not stored as a number)
$ws.Range("AE2").Value = "'2019"
